$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three pairs of data rows had their observation content swapped between
# each other (the row number / row-level formatting stays put, only the
# per-observation field values move). Only the columns whose values
# actually differ between the two rows of a pair are touched, so that
# unrelated (and especially blank placeholder) cells are left alone.

function Swap-Cell($sheet, $col, $r1, $r2) {
    $c1 = $sheet.Range($col + $r1)
    $c2 = $sheet.Range($col + $r2)
    $v1 = $c1.Value2
    $v2 = $c2.Value2
    $c1.Value2 = $v2
    $c2.Value2 = $v1
}

function Swap-Rows($sheet, $r1, $r2) {
    foreach ($col in @("A","B","E","F","G","H","Q","R")) {
        Swap-Cell $sheet $col $r1 $r2
    }

    # AJ/AK/AO ("Substratnamn" / "Vetenskapligt Substratnamn" /
    # "Substrat-beskrivning") move from one row of the pair to the other.
    $aj1 = $sheet.Range("AJ" + $r1)
    $ak1 = $sheet.Range("AK" + $r1)
    $ao1 = $sheet.Range("AO" + $r1)
    $aj2 = $sheet.Range("AJ" + $r2)
    $ak2 = $sheet.Range("AK" + $r2)
    $ao2 = $sheet.Range("AO" + $r2)

    $aj1v = $aj1.Value2
    $ak1v = $ak1.Value2
    $ao1v = $ao1.Value2
    $aj2v = $aj2.Value2
    $ak2v = $ak2.Value2
    $ao2v = $ao2.Value2

    if ($aj2v -eq $null) { $aj1.ClearContents() } else { $aj1.Value2 = $aj2v }
    if ($ak2v -eq $null) { $ak1.ClearContents() } else { $ak1.Value2 = $ak2v }
    if ($ao2v -eq $null) { $ao1.ClearContents() } else { $ao1.Value2 = $ao2v }

    if ($aj1v -eq $null) { $aj2.ClearContents() } else { $aj2.Value2 = $aj1v }
    if ($ak1v -eq $null) { $ak2.ClearContents() } else { $ak2.Value2 = $ak1v }
    if ($ao1v -eq $null) { $ao2.ClearContents() } else { $ao2.Value2 = $ao1v }
}

Swap-Rows $ws 16 17
Swap-Rows $ws 19 20
Swap-Rows $ws 24 25
